$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-coerced to a number by Excel's
# type inference (e.g. "18.80" -> 18.8, "0.000007496" -> 7.496E-06). Mark them
# as Text first so the literal string (incl. trailing zeros) round-trips, then
# restore the default "Normal" style so no stray formatting is left behind.
$textRefs = @("D2", "D3", "D5", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D32", "D34", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D45", "D46", "D47", "D48", "D49", "D50")
foreach ($ref in $textRefs) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = "30.543.57"
$ws.Range("E2").Value = "  -0.29%  "
$ws.Range("D3").Value = "1.874.86"
$ws.Range("E3").Value = "  -0.84%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "247.86"
$ws.Range("E5").Value = "  +0.89%  "
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("D7").Value = "0.4762"
$ws.Range("E7").Value = "  -0.56%  "
$ws.Range("D8").Value = "0.2908"
$ws.Range("E8").Value = "  +0.29%  "
$ws.Range("D9").Value = "0.06499"
$ws.Range("E9").Value = "  -1.03%  "
$ws.Range("D10").Value = "21.87"
$ws.Range("E10").Value = "  +1.15%  "
$ws.Range("D11").Value = "0.07754"
$ws.Range("E11").Value = "  -0.36%  "
$ws.Range("D12").Value = "0.7394"
$ws.Range("E12").Value = "  -0.20%  "
$ws.Range("D13").Value = "96.15"
$ws.Range("E13").Value = "  -1.42%  "
$ws.Range("D14").Value = "1.872.76"
$ws.Range("E14").Value = "  -1.41%  "
$ws.Range("D15").Value = "5.169"
$ws.Range("E15").Value = "  -0.28%  "
$ws.Range("D16").Value = "274.58"
$ws.Range("E16").Value = "  -2.13%  "
$ws.Range("D17").Value = "30.600.96"
$ws.Range("E17").Value = "  -0.18%  "
$ws.Range("D18").Value = "13.21"
$ws.Range("E18").Value = "  -1.91%  "
$ws.Range("D19").Value = "0.9998"
$ws.Range("E19").Value = "  -0.03%  "
$ws.Range("D20").Value = "0.000007496"
$ws.Range("E20").Value = "  -1.47%  "
$ws.Range("D21").Value = "2.118.22"
$ws.Range("E21").Value = "  -0.77%  "
$ws.Range("E22").Value = "  -0.10%  "
$ws.Range("D23").Value = "5.224"
$ws.Range("E23").Value = "  -1.17%  "
$ws.Range("D24").Value = "6.173"
$ws.Range("E24").Value = "  -0.82%  "
$ws.Range("D25").Value = "9.187"
$ws.Range("E25").Value = "  -1.78%  "
$ws.Range("D26").Value = "164.92"
$ws.Range("E26").Value = "  -0.59%  "
$ws.Range("D27").Value = "18.80"
$ws.Range("E27").Value = "  -1.65%  "
$ws.Range("D28").Value = "1.910"
$ws.Range("E28").Value = "  -2.93%  "
$ws.Range("D29").Value = "0.09849"
$ws.Range("E29").Value = "  -1.14%  "
$ws.Range("D30").Value = "1.336"
$ws.Range("E30").Value = "  -2.73%  "
$ws.Range("E31").Value = "  -1.38%  "
$ws.Range("D32").Value = "4.255"
$ws.Range("E32").Value = "  -2.31%  "
$ws.Range("E33").Value = "  -1.01%  "
$ws.Range("D34").Value = "0.04796"
$ws.Range("E34").Value = "  +0.41%  "
$ws.Range("E35").Value = "  -0.83%  "
$ws.Range("D36").Value = "0.6945"
$ws.Range("E36").Value = "  -1.03%  "
$ws.Range("D37").Value = "2.718"
$ws.Range("E37").Value = "  -0.07%  "
$ws.Range("D38").Value = "0.01855"
$ws.Range("E38").Value = "  -1.06%  "
$ws.Range("D39").Value = "2.758"
$ws.Range("E39").Value = "  -0.18%  "
$ws.Range("D40").Value = "6.284"
$ws.Range("E40").Value = "  -1.67%  "
$ws.Range("D41").Value = "73.26"
$ws.Range("E41").Value = "  +4.03%  "
$ws.Range("D42").Value = "1.979"
$ws.Range("E42").Value = "  +2.54%  "
$ws.Range("D43").Value = "0.4206"
$ws.Range("E43").Value = "  -0.26%  "
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("D45").Value = "0.8347"
$ws.Range("E45").Value = "  -1.41%  "
$ws.Range("D46").Value = "101.68"
$ws.Range("E46").Value = "  -0.77%  "
$ws.Range("D47").Value = "9.417"
$ws.Range("E47").Value = "  +0.84%  "
$ws.Range("B48").Value = "Elrond"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D48").Value = "35.29"
$ws.Range("E48").Value = "  -0.08%  "
$ws.Range("B49").Value = "Aptos"
$ws.Range("C49").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D49").Value = "6.971"
$ws.Range("E49").Value = "  -2.62%  "
$ws.Range("D50").Value = "913.16"
$ws.Range("E50").Value = "  -2.10%  "
$ws.Range("E51").Value = "  +0.93%  "

foreach ($ref in $textRefs) {
    $ws.Range($ref).Style = "Normal"
}

